$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("aglomerados")

# Carry the date-cell formatting from the previous row (A41) down to the
# new row's date cell (A42) before writing values, so the new row matches
# the existing table's look (same number format / style as other rows).
$ws.Range("A41").Copy()
$ws.Range("A42").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Append the new day's data as row 42
$ws.Cells.Item(42, 1).Value = 43982
$ws.Cells.Item(42, 2).Value = 2611
$ws.Cells.Item(42, 3).Value = 2485
$ws.Cells.Item(42, 4).Value = 41
$ws.Cells.Item(42, 5).Value = 159
$ws.Cells.Item(42, 6).Value = 696
$ws.Cells.Item(42, 7).Value = 63
$ws.Cells.Item(42, 8).Value = 248
$ws.Cells.Item(42, 9).Value = 0.179
$ws.Cells.Item(42, 10).Value = 847
$ws.Cells.Item(42, 11).Value = 0.064
$ws.Cells.Item(42, 12).Value = 81
$ws.Cells.Item(42, 13).Value = 78

# Restore the view: scroll the window so column D / row 28 is the top-left
# visible cell, then select N42 (mirroring the saved workbook's view state)
$excel.ActiveWindow.ScrollRow = 28
$excel.ActiveWindow.ScrollColumn = 4
$ws.Range("N42").Select()
